# Adds a new "2022" column (S) to the data table on the active sheet,
# mirroring the existing "2021" column (R) formatting, and fills in the
# new figures for row 3 (year header) through row 13 (data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: thin bottom-border spacer cell above the new year column ---
$ws.Range("S2").Value = $ws.Range("R2").Value
$ws.Range("S2").Borders.Item(9).LineStyle = $ws.Range("R2").Borders.Item(9).LineStyle
$ws.Range("S2").Borders.Item(9).Weight = $ws.Range("R2").Borders.Item(9).Weight
$ws.Range("S2").Font.Name = $ws.Range("R2").Font.Name
$ws.Range("S2").Font.Size = $ws.Range("R2").Font.Size
$ws.Range("S2").Font.Bold = $ws.Range("R2").Font.Bold

# --- Row 3: year header "2022" ---
$ws.Range("S3").Value = 2022
$ws.Range("S3").Font.Name = "Times New Roman"
$ws.Range("S3").Font.Size = 9
$ws.Range("S3").Font.Bold = $true
$ws.Range("S3").HorizontalAlignment = -4152
$ws.Range("S3").VerticalAlignment = -4108
$ws.Range("S3").Borders.Item(9).LineStyle = 1
$ws.Range("S3").Borders.Item(9).Weight = -4138

# --- Row 4: headline indicator value (bold) ---
$ws.Range("S4").Value = 17.700522048199787
$ws.Range("S4").NumberFormat = "0.0"
$ws.Range("S4").Font.Name = "Times New Roman"
$ws.Range("S4").Font.Size = 9
$ws.Range("S4").Font.Bold = $true
$ws.Range("S4").VerticalAlignment = -4108

# --- Rows 5-12: regular data rows ---
$ws.Range("S5").Value = 1.7610202290451711
$ws.Range("S6").Value = 3.9589300291403076
$ws.Range("S7").Value = 1.4859750619980623
$ws.Range("S8").Value = 1.1943569362276563
$ws.Range("S9").Value = 3.9154905266043296
$ws.Range("S10").Value = 0.84000241999604885
$ws.Range("S11").Value = 2.1393883316621789
$ws.Range("S12").Value = 1.8762854436950933

$dataRows = $ws.Range("S5:S12")
$dataRows.NumberFormat = "0.0"
$dataRows.Font.Name = "Times New Roman"
$dataRows.Font.Size = 9
$dataRows.Font.Bold = $false
$dataRows.VerticalAlignment = -4108

# --- Row 13: bottom (totals) row, medium bottom border ---
$ws.Range("S13").Value = 0.52907306983093583
$ws.Range("S13").NumberFormat = "0.0"
$ws.Range("S13").Font.Name = "Times New Roman"
$ws.Range("S13").Font.Size = 9
$ws.Range("S13").Font.Bold = $false
$ws.Range("S13").VerticalAlignment = -4108
$ws.Range("S13").Borders.Item(9).LineStyle = 1
$ws.Range("S13").Borders.Item(9).Weight = -4138

# Put the selection back on A1 (closest reachable approximation of the
# published sheet, which opens with no stale multi-cell selection).
$ws.Range("A1").Select()

Write-Host "done"
